$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# 1) Insert a new item row at row 7 for "AVIVAVASC 10/160MG 28 F.C. TAB."
#    (this pushes the existing items, total and footer down by one row)
# ---------------------------------------------------------------------
$ws.Rows("7:7").Insert()

# Clone the formatting (styles + row height + merges) of the row below
# (the former row 7, now shifted to row 8) onto the new blank row 7.
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)
$ws.Rows("7:7").RowHeight = $ws.Rows("8:8").RowHeight
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# Fill in the values for the new item.
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "AVIVAVASC 10/160MG 28 F.C. TAB."
Set-TextValue $ws.Range("H7") "0:3"
Set-TextValue $ws.Range("L7") "1"
Set-TextValue $ws.Range("N7") "140.00"
Set-TextValue $ws.Range("P7") "35.0000"
Set-TextValue $ws.Range("Q7") "0:1"

# Renumber the items that used to be #1-#5 (DEPOVIT ... سرنجات), which are
# now sitting in rows 8-12, to #2-#6.
for ($r = 8; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value() + 1
}

# ---------------------------------------------------------------------
# 2) Insert a new item row at row 13 (just above the totals row) for the
#    new "قياس سكر" item - pushes the totals & footer rows down by one.
# ---------------------------------------------------------------------
$ws.Rows("13:13").Insert()

$ws.Range("A12:Q12").Copy()
$ws.Range("A13:Q13").PasteSpecial(-4122)
$ws.Rows("13:13").RowHeight = $ws.Rows("12:12").RowHeight
$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "قياس سكر"
Set-TextValue $ws.Range("H13") "0:0"
Set-TextValue $ws.Range("L13") "0"
Set-TextValue $ws.Range("N13") "10.00"
Set-TextValue $ws.Range("P13") "10.0000"
Set-TextValue $ws.Range("Q13") "1:0"

# ---------------------------------------------------------------------
# 3) Update the totals row (now row 14) and the footer timestamp (row 15)
# ---------------------------------------------------------------------
$ws.Range("P14").Value = 158
$ws.Range("A15").Value = "Thursday, 28 August, 2025 10:01 AM"
